# The contents of rows 3-8 get cyclically permuted (a single 6-cycle):
#   new row 3 <- old row 4
#   new row 4 <- old row 5
#   new row 5 <- old row 6
#   new row 6 <- old row 8
#   new row 7 <- old row 3
#   new row 8 <- old row 7
#
# Rather than moving whole rows (which would risk Excel "helpfully"
# reinterpreting text such as the date strings in columns Y/Z/AA/AB as
# real dates), only the cells whose value actually differs between the
# source and destination row are touched here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (<- old row 4) ---
$ws.Range("A3").Value2 = 108369210
$ws.Range("Q3").Value2 = 663162.7064135609
$ws.Range("R3").Value2 = 6705981.337152475
$ws.Range("S3").Value2 = 25
$ws.Range("AC3").Value2 = "#SAKNAS!"
$ws.Range("AW3").Value2 = "Nadja Nilsson"
$ws.Range("AX3").Value2 = "Nadja Nilsson"

# --- Row 4 (<- old row 5) ---
$ws.Range("A4").Value2 = 108367419
$ws.Range("Q4").Value2 = 663092.6272863077
$ws.Range("R4").Value2 = 6705966.322238538
$ws.Range("S4").Value2 = 15
$ws.Range("AC4").ClearContents()
$ws.Range("AW4").Value2 = "Isac Carlsson"
$ws.Range("AX4").Value2 = "Isac Carlsson"

# --- Row 5 (<- old row 6) ---
$ws.Range("A5").Value2 = 108369211
$ws.Range("Q5").Value2 = 663155.2423936725
$ws.Range("R5").Value2 = 6706004.217932139
$ws.Range("S5").Value2 = 25
$ws.Range("AC5").Value2 = "#SAKNAS!"
$ws.Range("AW5").Value2 = "Nadja Nilsson"
$ws.Range("AX5").Value2 = "Nadja Nilsson"

# --- Row 6 (<- old row 8) ---
$ws.Range("A6").Value2 = 108367415
$ws.Range("B6").Value2 = 95511
$ws.Range("E6").Value2 = 221944
$ws.Range("F6").Value2 = "Lopplummer"
$ws.Range("G6").Value2 = "Huperzia selago"
$ws.Range("H6").Value2 = "(L.) Bernh. ex Schrank & Mart."
$ws.Range("Q6").Value2 = 663125.9516252303
$ws.Range("R6").Value2 = 6705852.723180643
$ws.Range("S6").Value2 = 15
$ws.Range("AC6").ClearContents()
$ws.Range("AW6").Value2 = "Isac Carlsson"
$ws.Range("AX6").Value2 = "Isac Carlsson"

# --- Row 7 (<- old row 3) ---
$ws.Range("A7").Value2 = 108367417
$ws.Range("Q7").Value2 = 663132.0536545257
$ws.Range("R7").Value2 = 6705596.601038971

# --- Row 8 (<- old row 7) ---
$ws.Range("A8").Value2 = 108367421
$ws.Range("B8").Value2 = 98520
$ws.Range("E8").Value2 = 222498
$ws.Range("F8").Value2 = "Blåsippa"
$ws.Range("G8").Value2 = "Hepatica nobilis"
$ws.Range("H8").Value2 = "Schreb."
$ws.Range("Q8").Value2 = 663073.8261583259
$ws.Range("R8").Value2 = 6705549.019092269
